# Apply the minor spec changes described in the commit:
#  - Remove the obsolete "T2.1 / X" placeholder row (old row 5)
#  - Append a trailing period to the T4.3 "circuit" specification text
#  - Adjust column A width now that the ID text is longer ("T2.1" row removed,
#    "T2.2" became "T2.1" etc. - the bestFit narrow width no longer applies)
#  - Refresh the active selection / scroll position to match the edited sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole row that held the placeholder "T2.1" / "X" / "X" entry.
# Everything below shifts up by one row (T4.3 text row -> row 13, etc.)
$deletedRowHeight = $ws.Rows.Item(5).RowHeight
$ws.Rows.Item(5).Delete()

# The three circuit-diagram pictures are anchored a few rows further down the
# sheet; keep them glued to the same cells they were on by moving them up by
# the height of the row that just disappeared (Excel does this automatically
# for "Move and size with cells" pictures when a row above them is deleted).
foreach ($shp in $ws.Shapes) {
    $shp.Top = $shp.Top - $deletedRowHeight
}

# Renumber the old "T2.2" row (now row 5) to "T2.1" since the real T2.1 row was removed.
$ws.Range("A5").Value = "T2.1"

# The footnote anchor name pointed at A15; after the row deletion that same
# physical cell is now A14, so re-point the defined name accordingly (Excel
# does this automatically when rows shift above a referenced cell).
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Blad1!_ftn1") {
        $n.RefersTo = "=Blad1!`$A`$14"
    }
}

# Correct the T4.3 specification text: add the missing trailing period.
$ws.Range("C13").Value = "De voedingsheaders leveren de directe voeding door aan de uitgangspoorten. Een minimale vorm van bescherming in de vorm van een diode en zekering wordt toegepast volgens onderstaande circuit."

# Column A no longer needs the narrow "best fit" width now that the ID texts
# changed; give it a plain custom width of 11 characters (stored column width
# is ColumnWidth + 5/6, so ask for 10 5/6 to land exactly on 11 in the XML).
$ws.Columns.Item(1).ColumnWidth = 10 + 1/6

# Update the selection to reflect the new extent of the table and reset the
# scrolled view back to the top-left (the sheet no longer needs to stay
# scrolled down to row 8).
$ws.Range("A1:D18").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
